# Quarterly indexing esoteric bug-fix operation
#
# Each data row (2..16) holds a rolling window of up-to-10 quarter-over-quarter
# error values (Q0 .. Q9, columns B..K). The series was being written one
# quarter out of phase: the newest observation needs to be *prepended* in
# column B, with the rest of the history shifting one column to the right
# and the oldest observation (once the window is full, i.e. it would land
# past column K) falling off the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlToLeft = -4159

# The newest Q0 error for every row, keyed by row number, in the same
# top-to-bottom order as the sheet.
$newValues = @{
    2  = -0.5825945370336409
    3  = 0.09567504080935779
    4  = -0.2604190369987228
    5  = 0.8354549961584912
    6  = -0.1000793599026215
    7  = -0.3537865060796963
    8  = 0.1481773904324453
    9  = 0.157445989004155
    10 = -0.5006594565260708
    11 = 0.2803578805354692
    12 = -0.1719748578450117
    13 = 0.3058625397463315
    14 = -0.6123299526872862
    15 = 0.6883713851991116
    16 = -0.2766911554241067
}

$firstDataCol = 2   # column B
$lastWindowCol = 11 # column K (Q0..Q9 == 10 columns wide)

foreach ($row in 2..16) {
    # Find the right-most populated column in this row (column A holds the
    # row label, so the data never starts before column B).
    $lastCell = $ws.Cells.Item($row, 200).End($xlToLeft)
    $lastCol = $lastCell.Column
    if ($lastCol -lt $firstDataCol) {
        $lastCol = $firstDataCol - 1
    }

    # Shift existing values one column to the right, working from the end
    # of the row backwards so values aren't clobbered before they're read.
    # Anything that would be pushed past the end of the fixed-width window
    # (column K) simply drops off.
    $shiftedLastCol = [Math]::Min($lastCol + 1, $lastWindowCol)
    for ($col = $shiftedLastCol; $col -gt $firstDataCol; $col--) {
        $srcCol = $col - 1
        if ($srcCol -ge $firstDataCol) {
            $srcValue = $ws.Cells.Item($row, $srcCol).Value2
            $ws.Cells.Item($row, $col).Value = $srcValue
        }
    }

    # Prepend the newest observation.
    $ws.Cells.Item($row, $firstDataCol).Value = $newValues[$row]
}
